# Bugfixed the naive forecaster component module
# The forecast vector table had an off-by-one: the first data row (old A2:E2,
# an AR(2) warm-up row with no y_0_forecast/y_1_forecast) is dropped and the
# remaining rows shift up by one. The y_1_forecast (column E) values are then
# recomputed for the shifted table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete leading row; everything below shifts up one row,
# shrinking the table from 19 data-ish rows (A1:E19) to 18 (A1:E18).
$ws.Rows.Item(2).Delete()

# Rows that now have fewer than 2 prior AR(2) observations no longer carry a
# y_1_forecast value.
$ws.Range("E2:E5").ClearContents()

# Recomputed y_1_forecast (column E) values for the remaining rows.
$ws.Range("E6").Value = 2.693188401769642
$ws.Range("E7").Value = 3.947916604971446
$ws.Range("E8").Value = 4.998814576944932
$ws.Range("E9").Value = 4.673582741620552
$ws.Range("E10").Value = 4.372458986620376
$ws.Range("E11").Value = 4.927320050172312
$ws.Range("E12").Value = 4.339089271348406
$ws.Range("E13").Value = 3.243024666552685
$ws.Range("E14").Value = 0.2915162802050064
$ws.Range("E15").Value = 3.818597641626909
$ws.Range("E16").Value = 0.6985632195332103
$ws.Range("E17").Value = 0.3452735157291054
$ws.Range("E18").Value = 1.5902148106679
